# Applies newly (re)generated experiment task orders to each sheet.
# Sheet tab order / rIds stay the same; only sheet names + the
# task_order lists they contain are refreshed with a new random draw.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1 (was GNG_TO...) -> TOL_TO..., gains two extra rows (4,5)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TOL_TO-16515890337749622"
$ws1.Range("B2").Value = "MM_stims-16515890337437139.csv"
$ws1.Range("B3").Value = "ZM_stims-16515890337280896.csv"
$ws1.Range("B4").Value = "MM_stims-16515890337593365.csv"
$ws1.Range("B5").Value = "ZM_stims-16515890337437139.csv"

# New rows 6 & 7 - copy formatting from row 5 so column A keeps style 1
$ws1.Range("A5:B5").Copy() | Out-Null
$ws1.Range("A6:B6").PasteSpecial(-4122) | Out-Null
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "MM_stims-16515890337749622.csv"

$ws1.Range("A5:B5").Copy() | Out-Null
$ws1.Range("A7:B7").PasteSpecial(-4122) | Out-Null
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "ZM_stims-16515890337593365.csv"

# ---------------------------------------------------------------
# Sheet 2 (NB_TO...) - new random order, same row count
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16515890352040193"
$ws2.Range("B2").Value = "OB-16515890347923858.csv"
$ws2.Range("B3").Value = "TB-16515890348705106.csv"
$ws2.Range("B4").Value = "OB-16515890347142923.csv"
$ws2.Range("B5").Value = "ZB-match_5-16515890339348633.csv"
$ws2.Range("B6").Value = "ZB-match_9-165158903414242.csv"
$ws2.Range("B7").Value = "TB-16515890351883948.csv"
$ws2.Range("B8").Value = "TB-16515890348392608.csv"
$ws2.Range("B9").Value = "ZB-match_2-16515890343996549.csv"
$ws2.Range("B10").Value = "OB-1651589034761171.csv"

# ---------------------------------------------------------------
# Sheet 3 (RS_TO...) - the two options swap places
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16515890352040193"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# ---------------------------------------------------------------
# Sheet 4 (was TOL_TO...) -> GNG_TO..., loses two rows (6,7)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "GNG_TO-16515890352352693"
$ws4.Range("B2").Value = "go_stims-16515890352040193.csv"
$ws4.Range("B3").Value = "GNG_stims-1651589035219645.csv"
$ws4.Range("B4").Value = "go_stims-1651589035219645.csv"
$ws4.Range("B5").Value = "GNG_stims-16515890352352693.csv"
$ws4.Range("A6:B7").Delete(-4162) | Out-Null

# ---------------------------------------------------------------
# Sheet 5 (vSAT_TO...) - new random order, same row count
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16515890352977805"
$ws5.Range("B2").Value = "SAT_stims-16515890352508934.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651589035266519.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651589035282145.csv"
$ws5.Range("B5").Value = "SAT_stims-16515890352352693.csv"
